$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the existing report table (pushes header/data
# from rows 1-2 down to rows 5-6) to make room for a title + export-time line.
$ws.Rows("1:4").Insert()

# Row 1: report title, merged across A1:B1, bold white text on an accent fill.
$ws.Range("A1").Value = "Báo cáo xuất kho thuốc"
$ws.Range("A1:B1").Merge()
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Size = 15
$ws.Range("A1:B1").Font.ThemeColor = 1
$ws.Range("A1:B1").Interior.ThemeColor = 8
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Rows(1).RowHeight = 19.5

# Row 3: "Thời gian:" label (bold) followed by the export-time placeholder.
$ws.Range("A3").Value = "Thời gian:"
$ws.Range("A3").Font.Bold = $true
$ws.Range("B3").Value = "&=ExportTime"

# AutoFilter over the (now shifted) header + data rows.
$ws.Range("A5:H6").AutoFilter()
$n = $ws.Names.Add('_xlnm._FilterDatabase', '=Sheet1!$A$5:$H$6')
$n.Visible = $false

# Restore normal view zoom and move the selection like the saved workbook.
$excel.ActiveWindow.Zoom = 100
$ws.Range("D5").Select()
